$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.701.95'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '2.230.73'
$ws.Range('E3').Value = '  -2.88%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = "'230.28"
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').Value = "'0.644"
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('D7').Value = "'63.30"
$ws.Range('E7').Value = '  +3.47%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = "'0.448"
$ws.Range('E9').Value = '  +5.19%  '
$ws.Range('D10').Value = "'0.0962"
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('D11').Value = "'56.79"
$ws.Range('E11').Value = '  -1.97%  '
$ws.Range('D12').Value = "'26.56"
$ws.Range('E12').Value = '  +8.94%  '
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').Value = '2.562.09'
$ws.Range('E14').Value = '  -2.78%  '
$ws.Range('D15').Value = "'15.39"
$ws.Range('E15').Value = '  -2.97%  '
$ws.Range('D16').Value = "'6.12"
$ws.Range('E16').Value = '  +3.91%  '
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = '2.227.99'
$ws.Range('E18').Value = '  -3.05%  '
$ws.Range('D19').Value = '43.586.86'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('D20').Value = '0.0₃0983'
$ws.Range('E20').Value = '  +3.73%  '
$ws.Range('D21').Value = "'72.59"
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('D22').Value = "'6.01"
$ws.Range('E22').Value = '  -4.15%  '
$ws.Range('D23').Value = "'248.51"
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = "'2.42"
$ws.Range('E25').Value = '  -5.61%  '
$ws.Range('D26').Value = "'3.40"
$ws.Range('E26').Value = '  +23.19%  '
$ws.Range('E27').Value = '  -2.72%  '
$ws.Range('D28').Value = "'9.84"
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('D29').Value = "'170.51"
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('D30').Value = "'20.78"
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('E31').Value = '  -2.24%  '
$ws.Range('E32').Value = '  -3.98%  '
$ws.Range('D33').Value = "'0.126"
$ws.Range('E33').Value = '  +2.71%  '
$ws.Range('D34').Value = "'0.0696"
$ws.Range('E34').Value = '  +5.33%  '
$ws.Range('D35').Value = "'4.73"
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('D36').Value = "'4.87"
$ws.Range('E36').Value = '  -4.37%  '
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('E38').Value = '  -2.71%  '
$ws.Range('D39').Value = "'2.26"
$ws.Range('E39').Value = '  -6.17%  '
$ws.Range('D40').Value = "'0.0258"
$ws.Range('E40').Value = '  +2.95%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = "'0.000219"
$ws.Range('E42').Value = '  -3.21%  '
$ws.Range('D43').Value = "'8.20"
$ws.Range('E43').Value = '  -7.77%  '
$ws.Range('D44').Value = "'17.08"
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').Value = "'96.93"
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = "'0.0942"
$ws.Range('E46').Value = '  -3.01%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = "'4.36"
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = "'1.18"
$ws.Range('E48').Value = '  -3.07%  '
$ws.Range('E49').Value = '  +2.88%  '
$ws.Range('D50').Value = '1.427.06'
$ws.Range('E50').Value = '  -3.76%  '
$ws.Range('E51').Value = '  +1.52%  '

# Reset style on forced-text numeric cells so no quote-prefix / custom number format lingers
$numericRefs = @('D5','D6','D7','D9','D10','D11','D12','D15','D16','D21','D22','D23','D25','D26','D28','D29','D30','D33','D34','D35','D36','D39','D40','D42','D43','D44','D45','D46','D47','D48')
foreach ($ref in $numericRefs) {
  $ws.Range($ref).Style = "Normal"
}
